$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Insert a new leave-card row for "10/2023" right after the existing
#    9/2023 row (sheet row 371). This pushes every later monthly row down
#    by one (old row 371 -> new row 372, old row 511 -> new row 512).
# ---------------------------------------------------------------------------
$ws.Rows.Item(371).Insert()

# The freshly inserted row has no formatting of its own yet - copy it from
# the row directly below (the shifted-down former row 371), which carries
# the same "blank monthly row" look used throughout the table.
$ws.Range("A372:K372").Copy()
$ws.Range("A371:K371").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Grow Table1 so the table (and its autofilter) cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K512"))

# Re-assert the calculated-column formula text for the two rows the resize
# touches (Excel's structured-reference shorthand resolves back to this
# canonical form on save, matching every other row in the column).
$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G371").Formula = $earnedFormula
$ws.Range("G512").Formula = $earnedFormula

# ---------------------------------------------------------------------------
# 2. Fill in the new leave records.
#    Row 370 (9/2023) earns an extra SL(7-0-0) entry covering 9/30-10/10.
#    Row 371 (new, no PERIOD date) records SL(4-0-0) covering 10/17-20.
# ---------------------------------------------------------------------------
$ws.Range("B370").Value = "SL(7-0-0)"
$ws.Range("C370").Value = 1.25
$ws.Range("H370").Value = 7
$ws.Range("K370").Value = "9/30 - 10/10/2023"

$ws.Range("B371").Value = "SL(4-0-0)"
$ws.Range("H371").Value = 4
$ws.Range("K371").Value = "10/17-20/2023"

# ---------------------------------------------------------------------------
# 3. Back-fill the EARNED (VL) column for the preceding months that were
#    still blank.
# ---------------------------------------------------------------------------
$ws.Range("C365").Value = 1.25
$ws.Range("C366").Value = 1.25
$ws.Range("C367").Value = 1.25
$ws.Range("C368").Value = 1.25
$ws.Range("C369").Value = 1.25

# ---------------------------------------------------------------------------
# 4. Leave the selection where the user ended up editing.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E9").Select()
$ws.Range("K371").Select()
